# Add the new "ODI Batting Extra" worksheet right after "ODI Batting" and
# populate it with the match-extras data (4s/6s/MoM breakdown).

$wb = $excel.ActiveWorkbook

$odiBatting = $wb.Worksheets.Item("ODI Batting")
$ws = $wb.Worksheets.Add($null, $odiBatting)
$ws.Name = "ODI Batting Extra"

# Header row - bold, thin border, centered/top-aligned (matches the style
# used for the header rows on the other sheets in this workbook).
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Row 2 - match 4735. Every field is text except BATTING_POSITION, which is
# a real number. Numeric-looking text is entered with a leading "'" so it
# is stored as text rather than being reinterpreted as a number/percentage.
$ws.Range("A2").Value = "'4735"
$ws.Cells.Item(2, 2).Value = 7
$ws.Range("C2").Value = "'4"
$ws.Range("D2").Value = "'1"
$ws.Range("E2").Value = "'17.88%"
$ws.Range("F2").Value = "NO"

# Row 3 - match 4745. MATCH_CODE and MAN_OF_MATCH are populated; the rest
# are present but blank (text cells with an empty value).
$ws.Range("A3").Value = "'4745"
$ws.Range("B3").Value = "'"
$ws.Range("C3").Value = "'"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "'"
$ws.Range("F3").Value = "NO"

$wb.Worksheets.Item("Player Info").Activate()
